# Weekly refresh: insert one new pair of daily-price rows (Primera/Segunda)
# for Perejil @ Vega Monumental Concepcion, ahead of the existing history.
# Everything below row 185 shifts down by two rows; all of the shifted
# rows keep their original values untouched (Excel does that for us).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 186, pushing the old
# rows 186:233 down to 188:235 (and bumping the sheet dimension to R235).
$ws.Range("A186:R187").Insert()

# --- New row 186 (Calidad = Primera) ---
$ws.Cells.Item(186, 1).Value = 11
$ws.Cells.Item(186, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(186, 3).Value = "Bíobío"
$ws.Cells.Item(186, 4).Value = 45135
$ws.Cells.Item(186, 5).Value = 8
$ws.Cells.Item(186, 6).Value = 100112044
$ws.Cells.Item(186, 7).Value = "Perejil"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 200
$ws.Cells.Item(186, 11).Value = 600
$ws.Cells.Item(186, 12).Value = 700
$ws.Cells.Item(186, 13).Value = 650
$ws.Cells.Item(186, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(186, 15).Value = "Región de Ñuble"
$ws.Cells.Item(186, 16).Value = 650
$ws.Cells.Item(186, 17).Value = 1
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# --- New row 187 (Calidad = Segunda) ---
$ws.Cells.Item(187, 1).Value = 11
$ws.Cells.Item(187, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(187, 3).Value = "Bíobío"
$ws.Cells.Item(187, 4).Value = 45135
$ws.Cells.Item(187, 5).Value = 8
$ws.Cells.Item(187, 6).Value = 100112044
$ws.Cells.Item(187, 7).Value = "Perejil"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Segunda"
$ws.Cells.Item(187, 10).Value = 100
$ws.Cells.Item(187, 11).Value = 500
$ws.Cells.Item(187, 12).Value = 500
$ws.Cells.Item(187, 13).Value = 500
$ws.Cells.Item(187, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(187, 15).Value = "Región de Ñuble"
$ws.Cells.Item(187, 16).Value = 500
$ws.Cells.Item(187, 17).Value = 1
$ws.Cells.Item(187, 18).Value = "Hortaliza"
